# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets.
# Mapping of row -> new F value, based on the authoritative diff.
$updates = @{
    3  = 1405
    6  = 239
    10 = 138
    11 = 4703
    12 = 6962
    18 = 4168
    19 = 980
    21 = 69
    25 = 175
    34 = 596
    36 = 552
    41 = 217
}

$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
